$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: 2021 year data ---------------------------------------
# Copy the formatting of the last existing year-label cell (A6, which
# carries the bold/bordered/centered style) onto the new label cells so
# we reuse the workbook's existing style slot instead of creating a new
# one.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 104.2
$ws.Range("C7").Value = 101.6
$ws.Range("D7").Value = 117
$ws.Range("E7").Value = 99.40000000000001
$ws.Range("F7").Value = 105.3
$ws.Range("G7").Value = 102.5
$ws.Range("H7").Value = 101

# --- New row 8: 2022 year data (partial - only B8 populated so far) --
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 105.2

# C8:H8 are present but still blank for 2022 (not yet reported) - give
# them an explicit empty-text value so the cells exist in the sheet
# (matching the other rows' shape) instead of being entirely absent,
# then strip the quote-prefix formatting that a leading apostrophe
# implies so the cells keep the workbook's default (unstyled) look.
$ws.Range("C8:H8").Value = "'"
$ws.Range("C8:H8").Style = "Normal"

$wb.Save()
